$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.054.56'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '3.261.35'
$ws.Range('E3').Value = '  +1.35%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('E7').Value = '  +4.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0956'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.82%  '
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('D13').Value = '3.778.67'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.29'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.00'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').Value = '3.261.93'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.19%  '
$ws.Range('D19').Value = '56.932.02'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000108'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '293.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.82%  '
$ws.Range('E25').Value = '  -1.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '28.16'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.53%  '
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '39.92'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0485'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('E36').Value = '  +1.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.32'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.27'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.18%  '
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.285'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.66%  '
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').Value = '  +4.50%  '
$ws.Range('D49').Value = '2.152.91'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.98'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.76%  '
$ws.Range('E51').Value = '  -5.59%  '
